# "all rows in data binding" -- the targetAssoc sheet had a stray/placeholder
# row (row 3: KRAS / assocCount 0) that shouldn't have been part of the data
# binding output. Remove it so the remaining rows shift up and the sheet's
# used range shrinks from A1:E5 to A1:E4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("targetAssoc")

# Delete entire row 3 - this shifts rows 4 and 5 up to become rows 3 and 4.
$ws.Rows.Item(3).Delete()

# Reflect the post-edit selection left behind in the file.
$ws.Range("C16").Select()
